$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.162.26"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.372.08"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'303.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'95.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'18.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "2.737.03"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "2.363.50"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "43.148.83"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "'68.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'235.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +15.10%  "
$ws.Range("D29").Value = "'9.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "'32.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'4.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").Value = "'2.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.07%  "
$ws.Range("D39").Value = "'123.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.52%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").Value = "'21.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.47%  "
$ws.Range("D43").Value = "1.937.64"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("E46").Value = "  -7.15%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "2.594.26"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D50").Value = "'71.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  +1.15%  "
